# Digital Transformation Summary - Portuguese copy-edit pass
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "Resumo da Transformação Digital" "Resumo da transformação digital"

Replace-Text "A Fabrikam Inc vem passando por uma iniciativa abrangente de transformação digital com o objetivo de melhorar a eficiência operacional, melhorar a experiência do cliente e impulsionar a inovação." `
             "A Fabrikam Inc. vem passando por uma iniciativa abrangente de transformação digital com o objetivo de aumentar a eficiência operacional, melhorar a experiência do cliente e promover a inovação."

Replace-Text "Implementação de Infraestrutura em Nuvem" "Implementação da infraestrutura de nuvem"

Replace-Text "Segurança de dados aprimorada e conformidade com os padrões do setor." "Maior segurança de dados e conformidade com os padrões do setor."

Replace-Text "Adoção de IA e Machine Learning" "Adoção de IA e aprendizado de máquina"

Replace-Text "Integração de análises alimentadas por IA para agilizar os processos de tomada de decisão." `
             "Integração de análises da plataforma de IA para agilizar os processos de tomada de decisão."

Replace-Text "Experiência Digital do Cliente" "Experiência digital do cliente"

Replace-Text "Automação de Processo" "Automação de processo"

Replace-Text "Implementação de automação robótica de processos (RPA) para tarefas rotineiras." `
             "Implementação da RPA (automação robótica de processos) para tarefas rotineiras."

Replace-Text "Redução de 40% no tempo de processamento das principais operações de negócios." `
             "Obteve uma redução de 40% no tempo de processamento das principais operações de negócios."

Replace-Text "Treinamento e Desenvolvimento de Funcionários" "Treinamento e desenvolvimento de funcionários"

Replace-Text "Conduzi programas de alfabetização digital para todos os funcionários." `
             "Conduziu programas de alfabetização digital para todos os funcionários."

Replace-Text "Aumento do engajamento dos funcionários e adoção de novas ferramentas em 35%." `
             "Aumento da participação do funcionário e adoção de novas ferramentas em 35%."

Replace-Text "1º trimestre de 2024: Migração concluída para infraestrutura em nuvem." `
             "1º trimestre de 2024: migração concluída para a infraestrutura de nuvem."

Replace-Text "2º trimestre de 2024: Lançada plataforma de análise orientada por IA." `
             "2º trimestre de 2024: lançamento da plataforma de análise orientada por IA."

Replace-Text "3º trimestre de 2024: Introduzido novo portal digital do cliente." `
             "3º trimestre de 2024: introdução do novo portal digital do cliente."

Replace-Text "4º trimestre de 2024: Atingi 50% de automação dos processos de rotina." `
             "4º trimestre de 2024: alcançou 50% de automação dos processos de rotina."

Replace-Text "Continue expandindo aplicativos de IA e aprendizado de máquina em todos os departamentos." `
             "Continue expandindo os aplicativos de IA e de aprendizado de máquina em todos os departamentos."

Replace-Text "Aprimore ainda mais a experiência digital do cliente com novos recursos e serviços." `
             "Melhore ainda mais a experiência digital do cliente com novos recursos e serviços."

Replace-Text "Concentre-se em medidas de segurança cibernética para proteger contra ameaças em evolução." `
             "Foque em medidas de segurança cibernética para proteção contra ameaças em evolução."

Replace-Text "Desenvolver uma estratégia digital abrangente para os próximos cinco anos." `
             "Desenvolva uma estratégia digital abrangente para os próximos cinco anos."

Replace-Text "A jornada de transformação digital da Fabrikam Inc." `
             "A jornada de transformação digital da Fabrikam Inc. levou a melhorias significativas na eficiência, na satisfação do cliente e no desempenho geral dos negócios."

Replace-Text "A organização continua comprometida em alavancar a tecnologia para impulsionar o crescimento e a inovação futuros." `
             "A organização continua comprometida em utilizar a tecnologia para impulsionar o crescimento e a inovação futuros."

# Bold the two Heading3/Heading2 runs that were previously non-bold
foreach ($para in $d.Paragraphs) {
    $text = $para.Range.Text.Trim()
    if ($text -eq "Automação de processo" -or $text -eq "Planos futuros") {
        $para.Range.Font.Bold = 1
    }
}
